$wb = $excel.ActiveWorkbook

# The "想去人数" (number of people wanting to go) counts were refreshed for
# three events, on both the "展览" sheet and the aggregated "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1450
    $ws.Range("F4").Value = 97
    $ws.Range("F6").Value = 21
}
